# Updated cryptos list on Sun Jun 23 18:39:29 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.079.86'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '3.474.23'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.98'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.67'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +5.65%  '
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = '4.064.18'
$ws.Range('E12').Value = '  -0.69%  '
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000178'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.22%  '
$ws.Range('D15').Value = '3.474.05'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').Value = '64.053.95'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.13'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.99'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.98%  '
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.42'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '385.11'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.568'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').Value = '3.613.40'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.56'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000111'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.97%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.22'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.12'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.80%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.97'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.63%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.43'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.35%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.153'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('B33').Value = 'RenzoRestakedETH'
$ws.Range('C33').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D33').Value = '3.501.28'
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '22.95'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.01%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.23'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.77'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.50'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.20%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '161.96'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.10%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0779'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.797'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.07%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.17'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.84%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.30'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.16%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.62'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.66'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.34%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.13'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.87%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.72'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.902'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.324.59'
$ws.Range('E50').Value = '  -5.52%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0254'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.70%  '
